# Update from MV -datos- : refresh rows 41-45 with revised quarterly figures
# and append a new row (46) for period 01-04-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 41-45 (columns B..X)
$updates = @{
    41 = @(118348, -568436, 301, 568737, -28, 0, 28, 0, 686812, 0, 0, 0, 0, 0, -161702, 0, 0, 0, -161702, 1527, 163230, 0, 280050)
    42 = @(-68225, 23494, 900, -22594, -6, 0, 6, 0, -91714, 0, 0, 0, 0, 0, -1323, 0, 0, 0, -1323, 3110, 4433, 0, -66902)
    43 = @(156324, 5990, 50, -5940, 18, 12, -6, 0, 150316, 0, 0, 0, 0, 0, 1482, 0, 0, 0, 1482, 2919, 1437, 0, 154842)
    44 = @(375, 6605, 1894, -4711, 1, 132, 130, 0, -6231, 0, 0, 0, 0, 0, 858, 0, 0, 0, 858, 4721, 3863, 0, -483)
    45 = @(250531, -654697, 504, 655200, -10, 0, 10, 0, 905237, 0, 0, 0, 0, 0, -162498, 0, 0, 0, -162498, 1195, 163694, 0, 413030)
}

$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $vals[$i]
    }
}

# New row 46: period 01-04-2021.
# Force text type (matching the "Serie" column's shared-string cells above it)
# instead of letting Excel auto-convert the dd-mm-yyyy-looking text to a date
# serial, then drop the number-format override so the cell keeps the sheet's
# default (unstyled) look, same as A2:A45.
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "01-04-2021"
$ws.Range("A46").ClearFormats()

$newRowVals = @(-67568, 11861, -72, -11933, 87, 88, 0, 0, -79516, 0, 0, 0, 0, 0, 368, 0, 0, 0, 368, 2627, 2259, 0, -67936)
for ($i = 0; $i -lt $columns.Length; $i++) {
    $ws.Range("$($columns[$i])46").Value = $newRowVals[$i]
}
